$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'20.447.61"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "'1.469.16"
$ws.Range("E3").Value = "  +4.53%  "
$ws.Range("E4").Value = "  +0.99%  "
$ws.Range("D5").Value = "'280.44"
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("D6").Value = "'0.8968"
$ws.Range("E6").Value = "  -10.52%  "
$ws.Range("D7").Value = "'0.3730"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "'0.3179"
$ws.Range("E8").Value = "  +3.65%  "
$ws.Range("D9").Value = "'39.51"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'1.050"
$ws.Range("D11").Value = "'0.06615"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "'5.536"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "'17.85"
$ws.Range("E14").Value = "  +5.70%  "
$ws.Range("D15").Value = "'6.190"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").Value = "'1.482.83"
$ws.Range("E16").Value = "  +5.44%  "
$ws.Range("D17").Value = "'0.00001031"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("D18").Value = "'0.05672"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").Value = "'0.9001"
$ws.Range("E19").Value = "  -10.23%  "
$ws.Range("D20").Value = "'70.39"
$ws.Range("E20").Value = "  -4.47%  "
$ws.Range("D21").Value = "'5.678"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "'14.62"
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("D23").Value = "'11.20"
$ws.Range("E23").Value = "  +3.23%  "
$ws.Range("D24").Value = "'2.289"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("D25").Value = "'20.609.52"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("D26").Value = "'2.261"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").Value = "'137.29"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  +3.19%  "
$ws.Range("D29").Value = "'1.646.42"
$ws.Range("E29").Value = "  +5.18%  "
$ws.Range("D30").Value = "'112.95"
$ws.Range("E30").Value = "  +3.66%  "
$ws.Range("D31").Value = "'3.935"
$ws.Range("E31").Value = "  +3.44%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.085"
$ws.Range("E32").Value = "  -5.67%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.8316"
$ws.Range("E33").Value = "  -1.63%  "
$ws.Range("D34").Value = "'0.07798"
$ws.Range("D35").Value = "'0.06065"
$ws.Range("E35").Value = "  +4.98%  "
$ws.Range("D36").Value = "'1.468"
$ws.Range("E36").Value = "  +15.97%  "
$ws.Range("D37").Value = "'4.833"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "'1.154"
$ws.Range("E38").Value = "  +8.93%  "
$ws.Range("D39").Value = "'10.48"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").Value = "'0.02034"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'0.1861"
$ws.Range("E41").Value = "  -3.04%  "
$ws.Range("D42").Value = "'0.9203"
$ws.Range("E42").Value = "  -8.15%  "
$ws.Range("D43").Value = "'0.5353"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").Value = "'3.578"
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'6.939"
$ws.Range("E45").Value = "  -17.49%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.21"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'122.74"
$ws.Range("E47").Value = "  +12.00%  "
$ws.Range("D48").Value = "'0.5252"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("D49").Value = "'1.816"
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").Value = "'0.06428"
$ws.Range("E50").Value = "  +4.41%  "
$ws.Range("D51").Value = "'1.030"
$ws.Range("E51").Value = "  -1.54%  "
